# Update the "as_of_utc" timestamp column (AA) for rows 2-26 on both the
# "Главные" and "Линейные" worksheets, changing the stamp from
# 2025-12-01 03:05:18 -> 2025-12-01 07:04:51.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-12-01 07:04:51"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
